# Add a new "segments" column (B) holding the segment names that used to
# live in column A, shift the two numeric columns (PercActivations,
# PercSegmentAreas) one slot to the right, and turn column A into a
# 0-based numeric index. This matches the commit "Added all current Data".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 20

# Step 1: shift column C (PercSegmentAreas) into the new column D.
# Range.Copy carries the value *and* the cell's existing style, so the
# header (s=1) and data cells (no explicit style) land correctly.
for ($r = 1; $r -le $lastRow; $r++) {
    $ws.Range("C$r").Copy($ws.Range("D$r"))
}

# Step 2: shift column B (PercActivations) into column C, same way.
for ($r = 1; $r -le $lastRow; $r++) {
    $ws.Range("B$r").Copy($ws.Range("C$r"))
}

# Step 3: build the new "segments" header in B1. Grab the header style by
# copying the (already relocated) C1 header cell, then overwrite its text.
$ws.Range("C1").Copy($ws.Range("B1"))
$ws.Range("B1").Value = "segments"

# Step 4: fill column B's data rows with the segment names that are still
# sitting in column A - plain assignment leaves these cells unstyled,
# exactly like the target.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Range("B$r").Value = $ws.Range("A$r").Value2
}

# Step 5: replace column A's text with a 0-based numeric row index while
# keeping the style (bold/border/center) it already had.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Range("A$r").Value = $r - 2
}
